# Generate Report for Handback
#
# For the "390884ee-b5f3-4a65-bcae-7f1ee6e7818d.md" source file (row 7 of
# both the zh-cn and the de-de sheets) the handback tool detected that the
# handback file it received is stale: a newer commit exists upstream. This
# records that detection in the report:
#   - I7 (Latest Target File)      -> the source file name, now a hyperlink
#   - J7 (Latest Handback File)    -> the (stale) generated xliff file name
#   - K7 (Latest Handback DateTime)-> the datetime the (stale) handback was generated
#   - P7 (Error Detail)            -> human readable "stale handback" message
# Column P ("Error Detail") is widened so the message is readable.

$wb = $excel.ActiveWorkbook

$sourceFile = "390884ee-b5f3-4a65-bcae-7f1ee6e7818d.md"
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80186c67218549afbe64c8d609cc2ca1dacafb62/e2e/390884ee-b5f3-4a65-bcae-7f1ee6e7818d.md"
$latestUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdea68905d1797bd23a8ffa9a21aeca5e3312e8a/e2e/390884ee-b5f3-4a65-bcae-7f1ee6e7818d.md"
$errorDetail = "The version of handback file is not the latest, current: $currentUrl, latest: $latestUrl."
$targetHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bdea68905d1797bd23a8ffa9a21aeca5e3312e8a/e2e/390884ee-b5f3-4a65-bcae-7f1ee6e7818d.md"

$sheets = @(
    @{ Name = "zh-cn"; Xlf = "390884ee-b5f3-4a65-bcae-7f1ee6e7818d.26b1871ff9d77e5c5098f5c04278a3a248c25a0d.zh-cn.xlf"; HandbackTime = "2016-09-05 08:59:17" },
    @{ Name = "de-de"; Xlf = "390884ee-b5f3-4a65-bcae-7f1ee6e7818d.26b1871ff9d77e5c5098f5c04278a3a248c25a0d.de-de.xlf"; HandbackTime = "2016-09-05 08:59:26" }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen the "Error Detail" column (P) so the message is legible.
    $ws.Columns.Item(16).ColumnWidth = 40

    # I7: Latest Target File -> hyperlink to the source .md file
    $ws.Hyperlinks.Add($ws.Range("I7"), $targetHyperlinkUrl, "", "", $sourceFile)

    # J7: Latest Handback File
    $ws.Range("J7").Value = $info.Xlf

    # K7: Latest Handback DateTime
    $ws.Range("K7").Value = $info.HandbackTime

    # P7: Error Detail
    $ws.Range("P7").Value = $errorDetail
}
